# Software BOM update: remove the extra (deployment/hosting-related)
# software entries - Amazon Web Services, Ubuntu Server, Apache HTTP
# Server, MariaDB and PHP - leaving only the original development-tool
# rows plus the SimplePoly City asset row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7-11 hold: Amazon Web Services, Ubuntu Server, Apache HTTP Server,
# MariaDB and PHP. Deleting them shifts the following row (SimplePoly
# City, previously row 12) up to become row 7.
$ws.Rows("7:11").Delete()

# Restore the selected cell to A2 (matches the saved workbook state).
$ws.Range("A2").Select()
